$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: move C1 value to D1 (wrapped in angle brackets), clear C1 (content + style)
$baseValue = $ws.Range("C1").Value2
$ws.Range("D1").Value = "<" + $baseValue + ">"
$ws.Range("C1").Clear()

# Row 2: B2 @data-namespace -> @prefix, C2 -> ":data", D2 -> "<http://sales.data/purchases#>"
$ws.Range("B2").Value = "@prefix"
$ws.Range("C2").Clear()
$ws.Range("C2").Value = ":data"
$ws.Range("D2").Value = "<http://sales.data/purchases#>"

# Row 3: B3 @schema-namespace -> @prefix, C3 -> ":schema", D3 -> "<http://sales.data/schema#>"
$ws.Range("B3").Value = "@prefix"
$ws.Range("C3").Clear()
$ws.Range("C3").Value = ":schema"
$ws.Range("D3").Value = "<http://sales.data/schema#>"

# Update selection to B1
$ws.Range("B1").Select()

# Add column E width (closest achievable value to target 9.140625 given
# the runtime's internal pixel-based quantization of ColumnWidth)
$ws.Columns.Item(5).ColumnWidth = 8.25

$wb.Save()
